# Apply "current state (new derivation)" value updates across the
# gens / lines / bus sheets (plain data cells, no formulas involved).
$wb = $excel.ActiveWorkbook
$wsGens = $wb.Worksheets.Item("gens")
$wsLines = $wb.Worksheets.Item("lines")
$wsBus = $wb.Worksheets.Item("bus")


# --- gens ---
$wsGens.Range("C11").Value = 100
$wsGens.Range("C24").Value = 289

# --- lines ---
$wsLines.Range("C2").Value = 20.950241
$wsLines.Range("D2").Value = 0
$wsLines.Range("C3").Value = -83.329948
$wsLines.Range("D3").Value = 0
$wsLines.Range("C4").Value = -45.620293
$wsLines.Range("D4").Value = 0
$wsLines.Range("C5").Value = -51.790415
$wsLines.Range("D5").Value = 0
$wsLines.Range("C6").Value = -24.259344
$wsLines.Range("D6").Value = 0
$wsLines.Range("C7").Value = -14.988768
$wsLines.Range("D7").Value = 0
$wsLines.Range("C8").Value = -248.34118
$wsLines.Range("D8").Value = -0
$wsLines.Range("C9").Value = -125.79041
$wsLines.Range("D9").Value = 0
$wsLines.Range("C10").Value = -116.62029
$wsLines.Range("D10").Value = 0
$wsLines.Range("C11").Value = -160.25934
$wsLines.Range("D11").Value = 0
$wsLines.Range("C12").Value = 175
$wsLines.Range("F12").Value = 0
$wsLines.Range("C13").Value = -13.836278
$wsLines.Range("D13").Value = 0
$wsLines.Range("C14").Value = 17.836278
$wsLines.Range("D14").Value = 0
$wsLines.Range("C15").Value = -145.2463
$wsLines.Range("D15").Value = 0
$wsLines.Range("C16").Value = -184.36917
$wsLines.Range("D16").Value = 0
$wsLines.Range("C17").Value = -207.46024
$wsLines.Range("D17").Value = 0
$wsLines.Range("C18").Value = -246.58312
$wsLines.Range("D18").Value = 0
$wsLines.Range("C19").Value = -269.74375
$wsLines.Range("D19").Value = 0
$wsLines.Range("C20").Value = -82.962788
$wsLines.Range("D20").Value = -0.000000000000034194518
$wsLines.Range("C21").Value = -201.27873
$wsLines.Range("D21").Value = 0
$wsLines.Range("C22").Value = -229.67355
$wsLines.Range("D22").Value = 0
$wsLines.Range("C23").Value = -145.02248
$wsLines.Range("D23").Value = 0
$wsLines.Range("C24").Value = -276.96279
$wsLines.Range("D24").Value = -0.000000000000048035156
$wsLines.Range("C25").Value = -2.0701978
$wsLines.Range("D25").Value = 0
$wsLines.Range("C26").Value = -204.13549
$wsLines.Range("D26").Value = 0.000000000000014210855
$wsLines.Range("C27").Value = -204.13549
$wsLines.Range("D27").Value = -0.000000000000014210855
$wsLines.Range("C28").Value = 248.34118
$wsLines.Range("D28").Value = 0
$wsLines.Range("C29").Value = -247.72902
$wsLines.Range("D29").Value = 0
$wsLines.Range("E29").Value = 0
$wsLines.Range("C30").Value = 23.696033
$wsLines.Range("D30").Value = -0.000000000000018725569
$wsLines.Range("C31").Value = -109.4256
$wsLines.Range("D31").Value = 0
$wsLines.Range("C32").Value = -138.30342
$wsLines.Range("D32").Value = 0
$wsLines.Range("C33").Value = -76.7128
$wsLines.Range("D33").Value = 0
$wsLines.Range("C34").Value = -76.7128
$wsLines.Range("D34").Value = -0
$wsLines.Range("C35").Value = -78.651984
$wsLines.Range("D35").Value = -0.000000000000018355353
$wsLines.Range("C36").Value = -78.651984
$wsLines.Range("D36").Value = -0.000000000000014210855
$wsLines.Range("C37").Value = -142.65198
$wsLines.Range("D37").Value = -0.0000000000000037005593
$wsLines.Range("C38").Value = -142.65198
$wsLines.Range("D38").Value = -0.000000000000014210855
$wsLines.Range("C39").Value = -161.69658
$wsLines.Range("D39").Value = 0
$wsLines.Range("C40").Value = 0
$wsLines.Range("D40").Value = -0.000000000000007327396600000001

# --- bus ---
$wsBus.Range("B2").Value = 100
$wsBus.Range("C2").Value = -57.131627
$wsBus.Range("B3").Value = 100
$wsBus.Range("C3").Value = -57.42493
$wsBus.Range("B4").Value = 100
$wsBus.Range("C4").Value = -39.549008
$wsBus.Range("B5").Value = 100
$wsBus.Range("C5").Value = -50.847548
$wsBus.Range("B6").Value = 100
$wsBus.Range("C6").Value = -53.253902
$wsBus.Range("B7").Value = 100
$wsBus.Range("C7").Value = -52.767136
$wsBus.Range("C8").Value = -29.37333
$wsBus.Range("B9").Value = 100
$wsBus.Range("C9").Value = -40.04833
$wsBus.Range("B10").Value = 100
$wsBus.Range("C10").Value = -37.765345
$wsBus.Range("B11").Value = 100
$wsBus.Range("C11").Value = -42.991316
$wsBus.Range("B12").Value = 100
$wsBus.Range("C12").Value = -25.564656
$wsBus.Range("B13").Value = 100
$wsBus.Range("C13").Value = -22.278335
$wsBus.Range("B14").Value = 100
$wsBus.Range("C14").Value = -12.616956
$wsBus.Range("B15").Value = 100
$wsBus.Range("C15").Value = -22.080219
$wsBus.Range("B16").Value = 100
$wsBus.Range("C16").Value = -5.7746076
$wsBus.Range("B17").Value = 100
$wsBus.Range("C17").Value = -5.7394142
$wsBus.Range("B18").Value = 100
$wsBus.Range("C18").Value = 0.70154024
$wsBus.Range("B19").Value = 100
$wsBus.Range("C19").Value = 2.2334986
$wsBus.Range("B20").Value = 100
$wsBus.Range("C20").Value = -6.284423
$wsBus.Range("B21").Value = 100
$wsBus.Range("C21").Value = -3.1383436
$wsBus.Range("B22").Value = 100
$wsBus.Range("C22").Value = 4.2280315
$wsBus.Range("B23").Value = 100
$wsBus.Range("C23").Value = 15.223399
$wsBus.Range("B24").Value = 100
$wsBus.Range("C24").Value = 0
$wsBus.Range("B25").Value = 100
$wsBus.Range("C25").Value = -18.688349
$wsBus.Range("B26").Value = 100
$wsBus.Range("C26").Value = 0
